$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in the sheet
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)   # D: total days
    $eCell = $ws.Cells.Item($r, 5)   # E: remaining days
    $fCell = $ws.Cells.Item($r, 6)   # F: start date (yyyymmdd as number)

    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    if ($dVal -eq $null -or $eVal -eq $null -or $fVal -eq $null) {
        continue
    }

    # Skip rows whose start date isn't a well-formed 8-digit yyyymmdd value
    $fStr = [string]([int64]$fVal)
    if ($fStr.Length -ne 8) {
        continue
    }

    if ([int]$eVal -eq 1) {
        # Cycle finished: reset remaining to total days and roll the start date forward by 10 days
        $eCell.Value2 = $dVal
        $fCell.Value2 = $fVal + 10
    } else {
        # One more day elapsed: decrement remaining days
        $eCell.Value2 = $eVal - 1
    }
}
